$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.942.33'
$ws.Range("D3").Value = '2.822.42'
$ws.Range("E3").Value = '  +2.87%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''355.54'
$ws.Range("E5").Value = '  +6.94%  '
$ws.Range("D6").Value = '''113.67'
$ws.Range("E6").Value = '  -1.93%  '
$ws.Range("E7").Value = '  +2.79%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +4.94%  '
$ws.Range("D10").Value = '''42.01'
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").Value = '''20.05'
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").Value = '''7.72'
$ws.Range("E14").Value = '  +1.45%  '
$ws.Range("D15").Value = '3.246.43'
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").Value = '2.829.69'
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").Value = '51.845.44'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = '''7.40'
$ws.Range("E19").Value = '  +8.18%  '
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").Value = '''13.57'
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("E22").Value = '  +2.35%  '
$ws.Range("D23").Value = '''269.98'
$ws.Range("E23").Value = '  -3.05%  '
$ws.Range("D24").Value = '''69.74'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").Value = '''2.79'
$ws.Range("E25").Value = '  +5.51%  '
$ws.Range("D26").Value = '''26.83'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("E29").Value = '  +1.60%  '
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").Value = '''50.82'
$ws.Range("E31").Value = '  +1.69%  '
$ws.Range("D32").Value = '''33.91'
$ws.Range("E32").Value = '  -3.21%  '
$ws.Range("E33").Value = '  +31.51%  '
$ws.Range("E34").Value = '  +5.56%  '
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  +0.61%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''4.90'
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''3.22'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").Value = '''18.41'
$ws.Range("E40").Value = '  -3.72%  '
$ws.Range("D41").Value = '''23.71'
$ws.Range("E41").Value = '  +2.99%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''2.58'
$ws.Range("E42").Value = '  +6.43%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''128.38'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("D46").Value = '''3.35'
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("D47").Value = '2.075.63'
$ws.Range("E47").Value = '  -0.54%  '
$ws.Range("E48").Value = '  +4.16%  '
$ws.Range("D49").Value = '''0.958'
$ws.Range("E49").Value = '  +9.82%  '
$ws.Range("D50").Value = '''5.70'
$ws.Range("E50").Value = '  +3.24%  '
$ws.Range("D51").Value = '''60.51'
$ws.Range("E51").Value = '  +1.11%  '
